$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily data point was recorded for 2026/02/04 before the existing
# rows continue with 2026/12/29 onward. Insert a new row at 756 which
# shifts the existing rows 756:797 down to 757:798.
$ws.Rows.Item(756).Insert()

# Populate the newly inserted row with the new data point.
# The leading apostrophe forces text entry for the date-like string so
# Excel does not auto-convert it into a date serial number; resetting
# the style afterwards removes the "quote prefix" formatting flag so the
# cell ends up as a plain, unstyled text cell like its neighbours.
$ws.Cells.Item(756, 1).Value = "'2026/02/04"
$ws.Cells.Item(756, 1).Style = "Normal"
$ws.Cells.Item(756, 2).Value = "水"
$ws.Cells.Item(756, 3).Value = 23
$ws.Cells.Item(756, 4).Value = 36
